# Append the two new daily GSC export rows to the "Chart" sheet.
# Row 75: 2025-12-18 (Invalid=0, Valid=31)
# Row 76: 2025-12-19 (Invalid=0, Valid=32)
#
# The Date column stores dates as plain text (shared strings), not real
# Excel dates, so we briefly mark the cell as Text ("@") before writing the
# literal string -- this stops Excel's automatic date-detection from
# converting it to a date serial number -- then clear the formatting again
# so the new cells end up with the same (default) style as the rest of the
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "2025-12-18"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = 0.0
$ws.Range("C75").Value = 31.0

$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "2025-12-19"
$ws.Range("A76").ClearFormats()
$ws.Range("B76").Value = 0.0
$ws.Range("C76").Value = 32.0
